$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B to fit the new Samsung/Apple test data
$ws.Columns.Item(1).ColumnWidth = 28.998697916666668
$ws.Columns.Item(2).ColumnWidth = 32.666666666666664

# Row 7: OPPO -> Samsung, and add the new Apple column value
$ws.Range("B7").Value = "Samsung"

# Row 8 (new): Samsung page verification message
$ws.Range("A8").Value = "samsungPageVerificationMessage"
$ws.Range("B8").Value = "Samsung Mobile Phones`n"
$ws.Range("B8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 32

# Row 9 (new): Apple mobile page verification message
$ws.Range("A9").Value = "AppleMobPageVerificationMessage"
$ws.Range("B9").Value = "Apple Store"

# Back to row 7 to add the new Apple product type (kept last so the
# "Apple" shared string is appended after the row 8/9 strings)
$ws.Range("C7").Value = "Apple"

# Match the author's final selection
$ws.Range("C7").Select()
